# This script rotates the content of rows 18, 19 and 20 (keeping row
# identity / row-local fields like "Osäker artbestämning" (AE) in place)
# according to the target diff:
#   new row18 <- old row20 data (Garnlav / Alectoria sarmentosa)
#   new row19 <- old row18 data (Tretåig hackspett / färska spår)
#   new row20 <- old row19 data (Tretåig hackspett / bobygge)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 18 (becomes the old "Garnlav" record, row 20 data) ----
$ws.Range("A18").Value = 131039579
$ws.Range("B18").Value = 79243
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("M18").Value = ""
$ws.Range("P18").Value = "Gotvad, Dlr"
$ws.Range("Q18").Value = 479079
$ws.Range("R18").Value = 6792475
$ws.Range("S18").Value = 50
$ws.Range("AC18").Value = "Rikligt till måttligt i en radie av ca 50 meter, synfältet"
$ws.Range("AE18").Value = $false

# ---- Row 19 (becomes the old "Tretåig hackspett / färska spår" record, row 18 data) ----
$ws.Range("A19").Value = 131041965
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("M19").Value = "färska spår"
$ws.Range("P19").Value = "Tandbergsvasseln, Dlr"
$ws.Range("Q19").Value = 479096
$ws.Range("R19").Value = 6792085
$ws.Range("S19").Value = 10
$ws.Range("AC19").Value = ""
$ws.Range("AE19").Value = $false

# ---- Row 20 (becomes the old "Tretåig hackspett / bobygge" record, row 19 data) ----
$ws.Range("A20").Value = 131039828
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "bobygge"
$ws.Range("P20").Value = "Gotvad, Dlr"
$ws.Range("Q20").Value = 479059
$ws.Range("R20").Value = 6792254
$ws.Range("S20").Value = 10
$ws.Range("AC20").Value = ""
$ws.Range("AE20").Value = $true
